{"js": "// The underlying OOXML diff for this revision is a pure re-serialization:\n// every hunk in the canonical XML diff reorders existing attributes /\n// namespace declarations (e.g. `w:val=\"single\" w:sz=\"6\"` -> `w:sz=\"6\"\n// w:val=\"single\"`) without adding, removing, or changing any element,\n// attribute value, text run, or document part. The document's visible\n// content, formatting, styles, numbering, header/footer text, and section\n// properties are unchanged.\n//\n// Office.js has no API surface for controlling raw XML attribute\n// serialization order, and doing so is not a user-visible document edit,\n// so there is nothing to change through the Word JavaScript API here.\n// We simply touch the body to confirm the context is usable and sync,\n// leaving the document content untouched.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying OOXML diff for this revision is a pure re-serialization:\n# every hunk in the canonical XML diff reorders existing attributes /\n# namespace declarations (e.g. `w:val=\"single\" w:sz=\"6\"` -> `w:sz=\"6\"\n# w:val=\"single\"`) without adding, removing, or changing any element,\n# attribute value, text run, or document part. The document's visible\n# content, formatting, styles, numbering, header/footer text, and section\n# properties are unchanged.\n#\n# The Word COM object model has no surface for controlling raw XML\n# attribute serialization order, and doing so is not a user-visible\n# document edit, so there is nothing to change through COM here.\n# We simply touch the document to confirm it is reachable, leaving the\n# document content untouched.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
